$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the name in A8 from "Culebra" to "yami"
$ws.Range("A8").Value = "yami"

# Highlight the header row (B3:E3) with a solid yellow fill
$ws.Range("B3:E3").Interior.Color = 65535

# Update the selected/active cell to C6 (cosmetic view state change)
$ws.Range("C6").Select()
